$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the date value ("10.01.2020" -> "10.01.20") for every data row
# in both the s1cDNADate (column A) and s2cDNADate (column D) fields.
# Temporarily force the cells to plain-text formatting so Excel does not
# reinterpret the value as a date serial number, then restore the
# (default) formatting so no stray styles are left behind.
$rngA = $ws.Range("A2:A29")
$rngD = $ws.Range("D2:D29")
$rngA.NumberFormat = "@"
$rngD.NumberFormat = "@"

for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 1).Value = "10.01.20"
    $ws.Cells.Item($r, 4).Value = "10.01.20"
}

$rngA.ClearFormats()
$rngD.ClearFormats()

# Update the last active selection to mirror the saved workbook state.
$ws.Range("G34").Select()

$wb.Save()
